$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear cells whose content was removed by the edit
$ws.Range("B17").ClearContents()
$ws.Range("C17").ClearContents()
$ws.Range("B22").ClearContents()
$ws.Range("C22").ClearContents()
$ws.Range("A23").ClearContents()

# Set final cell values (text)
$ws.Range('B1').Value = 'Ementa atual:'
$ws.Range('C1').Value = 'Ementa modificada (dados modificados em vermelho):'
$ws.Range('B2').Value = 'LOQ4072'
$ws.Range('C2').Value = 'LOQ4072'
$ws.Range('A3').Value = 'Nome:'
$ws.Range('B3').Value = ' Processos de Separação por Membranas'
$ws.Range('C3').Value = ' Processos de Separação por Membranas'
$ws.Range('A4').Value = 'Name:'
$ws.Range('B4').Value = 'Membrane Separation Processes'
$ws.Range('C4').Value = 'Membrane Separation Processes'
$ws.Range('A5').Value = 'Créditos-aula:'
$ws.Range('B5').Value = '4'
$ws.Range('C5').Value = '4'
$ws.Range('A6').Value = 'Créditos-trabalho'
$ws.Range('B6').Value = '0'
$ws.Range('C6').Value = '0'
$ws.Range('A7').Value = 'Carga horária:'
$ws.Range('B7').Value = '60 h'
$ws.Range('C7').Value = '60 h'
$ws.Range('A8').Value = 'Ativação:'
$ws.Range('B8').Value = '01/01/2018'
$ws.Range('C8').Value = '01/01/2018'
$ws.Range('A9').Value = 'Semestre ideal:'
$ws.Range('B9').Value = 'EQD-10,EQN-12'
$ws.Range('C9').Value = 'EQD-10,EQN-12'
$ws.Range('A10').Value = 'Objetivos:'
$ws.Range('B10').Value = '787307 - Luis Fernando Figueiredo Faria'
$ws.Range('C10').Value = '787307 - Luis Fernando Figueiredo Faria'
$ws.Range('A11').Value = 'Objectives:'
$ws.Range('B11').Value = 'Provide basic knowledge of the Principles of Membrane Separation. Provide general information to understand the technology involved in the different types of Membrane Separation Processes and to know the materials most used in the manufacture of membranes and their respective manufacturing processes. To study the phenomenological components involved in the mechanisms of transport through the membranes, and to identify the main advantages, disadvantages and applications of this type of separation process.'
$ws.Range('C11').Value = 'Provide basic knowledge of the Principles of Membrane Separation. Provide general information to understand the technology involved in the different types of Membrane Separation Processes and to know the materials most used in the manufacture of membranes and their respective manufacturing processes. To study the phenomenological components involved in the mechanisms of transport through the membranes, and to identify the main advantages, disadvantages and applications of this type of separation process.'
$ws.Range('A12').Value = 'Docentes responsáveis:'
$ws.Range('A13').Value = 'Programa resumido:'
$ws.Range('B13').Value = 'Semestral'
$ws.Range('C13').Value = 'Semestral'
$ws.Range('A14').Value = 'Short syllabus:'
$ws.Range('B14').Value = 'Classification of membrane processes and their applications. Technical preparation of different types of polymeric membranes. Types of modules and its main features. Theoretical fundamentals of the membranes synthesis  by the inversion of phases; Influence of the synthesis variables on the transport characteristics of the membranes. Presentation of the different types of membrane processes. Applications. Design for a specific application.'
$ws.Range('C14').Value = 'Classification of membrane processes and their applications. Technical preparation of different types of polymeric membranes. Types of modules and its main features. Theoretical fundamentals of the membranes synthesis  by the inversion of phases; Influence of the synthesis variables on the transport characteristics of the membranes. Presentation of the different types of membrane processes. Applications. Design for a specific application.'
$ws.Range('A15').Value = 'Programa:'
$ws.Range('B15').Value = '01/01/2018'
$ws.Range('C15').Value = '01/01/2018'
$ws.Range('A16').Value = 'Syllabus:'
$ws.Range('B16').Value = '1. Introduction: History and definition of membranes separation processes: comparison with classical separation processes. Membrane processes: classification and applications. 2. Membranes: Definition; materials used; rating on the structure and the type of proposed separation; Microporous membranes: manufacturing techniques and characteristics. Synthesis of membranes by the phase inversion technique: thermodynamic and kinetic aspects. 3. Fundamentals of Membrane Separation Processes: transport mechanisms in dense and microporous membranes.4. Membranes and Commercial Modules: membrane geometry; Techniques for manufacturing flat membranes, hollow fiber and tubular membranes. Modules: plate type, spiral, tubular and hollow fiber.5. Commercial Membrane Separation Processes. Main applications.'
$ws.Range('C16').Value = '1. Introduction: History and definition of membranes separation processes: comparison with classical separation processes. Membrane processes: classification and applications. 2. Membranes: Definition; materials used; rating on the structure and the type of proposed separation; Microporous membranes: manufacturing techniques and characteristics. Synthesis of membranes by the phase inversion technique: thermodynamic and kinetic aspects. 3. Fundamentals of Membrane Separation Processes: transport mechanisms in dense and microporous membranes.4. Membranes and Commercial Modules: membrane geometry; Techniques for manufacturing flat membranes, hollow fiber and tubular membranes. Modules: plate type, spiral, tubular and hollow fiber.5. Commercial Membrane Separation Processes. Main applications.'
$ws.Range('A17').Value = 'Avaliação:'
$ws.Range('A18').Value = 'Método:'
$ws.Range('B18').Value = '787307 - Luis Fernando Figueiredo Faria'
$ws.Range('C18').Value = '787307 - Luis Fernando Figueiredo Faria'
$ws.Range('A19').Value = 'Critério:'
$ws.Range('B19').Value = '-Provas escritas; -participação e conteúdo de trabalho e seminário;'
$ws.Range('C19').Value = '-Provas escritas; -participação e conteúdo de trabalho e seminário;'
$ws.Range('A20').Value = 'Norma de recuperação:'
$ws.Range('B20').Value = 'Média Final = (Prova1 + Prova2 + Nota de Trabalho) / 3Média final mínima de aprovação = 5,0'
$ws.Range('C20').Value = 'Média Final = (Prova1 + Prova2 + Nota de Trabalho) / 3Média final mínima de aprovação = 5,0'
$ws.Range('A21').Value = 'Bibliografia:'
$ws.Range('B21').Value = '(Prova escrita + Média Final)/2         Nota Final mínima para aprovação= 5,0'
$ws.Range('C21').Value = '(Prova escrita + Média Final)/2         Nota Final mínima para aprovação= 5,0'
$ws.Range('A22').Value = 'Requisitos:'
$ws.Range("B23").Value = "LOQ4085 -  Operações Unitárias I  (Requisito fraco)`n"
$ws.Range("C23").Value = "LOQ4085 -  Operações Unitárias I  (Requisito fraco)`n"

# Delete the now-unused last row (24), shrinking the sheet to A1:C23
$ws.Rows(24).Delete()

# Set row heights to match the target layout
$ws.Rows(13).RowHeight = 60
$ws.Rows(14).RowHeight = 60
$ws.Rows(15).RowHeight = 120
$ws.Rows(16).RowHeight = 120
$ws.Rows(18).RowHeight = 60
$ws.Rows(19).RowHeight = 60
$ws.Rows(20).RowHeight = 60
$ws.Rows(21).RowHeight = 120
$ws.Rows(23).RowHeight = 30
